# Updates the cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for the rows whose figures changed in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.643.25"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.598.58"
$ws.Range("E3").Value = "  +1.16%  "

# Row 5 - BNB
$ws.Range("D5").Value = "210.69"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  +1.64%  "

# Row 7 - USDC (Volume only)
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  +0.23%  "

# Row 9 - Cardano (Volume only)
$ws.Range("E9").Value = "  -1.15%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.90%  "

# Row 11 - TRON (Volume only)
$ws.Range("E11").Value = "  +0.54%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "1.823.80"
$ws.Range("E12").Value = "  +0.98%  "

# Row 13 - Wrapped Ether
$ws.Range("D13").Value = "1.594.28"
$ws.Range("E13").Value = "  +0.81%  "

# Row 14 - Polkadot (Volume only)
$ws.Range("E14").Value = "  -0.48%  "

# Row 15 - Polygon (Volume only)
$ws.Range("E15").Value = "  -1.12%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.76"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17 - Wrapped BTC
$ws.Range("D17").Value = "26.642.62"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18 - Shiba Inu (price uses U+2083 SUBSCRIPT THREE between the zeros)
$subscriptThree = [char]0x2083
$ws.Range("D18").Value = "0.0{0}0729" -f $subscriptThree
$ws.Range("E18").Value = "  +0.23%  "

# Row 19 - Bitcoin Cash
$ws.Range("D19").Value = "208.83"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20 - Dai (Volume only)
$ws.Range("E20").Value = "  -0.15%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  +1.26%  "

# Row 22 - Uniswap (Volume only)
$ws.Range("E22").Value = "  +0.22%  "

# Row 23 - Toncoin (Volume only)
$ws.Range("E23").Value = "  -2.61%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  +0.21%  "

# Row 25 - Monero
$ws.Range("D25").Value = "145.84"
$ws.Range("E25").Value = "  -0.24%  "

# Row 26 - BinanceUSD (Volume only)
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - Cosmos (Volume only)
$ws.Range("E27").Value = "  -2.63%  "

# Row 28 - Stellar (Volume only)
$ws.Range("E28").Value = "  +2.39%  "

# Row 29 - Ethereum Classic
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30 - Hedera (Volume only)
$ws.Range("E30").Value = "  +1.16%  "

# Row 31 - PancakeSwap (Volume only)
$ws.Range("E31").Value = "  +0.03%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "0.657"
$ws.Range("E33").Value = "  +0.33%  "

# Row 34 - Internet Computer (DFINITY) (Volume only)
$ws.Range("E34").Value = "  -0.11%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.295.31"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36 - Huobi Token (Volume only)
$ws.Range("E36").Value = "  +0.46%  "

# Row 37 - Lido DAO Token (Volume only)
$ws.Range("E37").Value = "  -1.40%  "

# Row 38 - VeChain (Volume only)
$ws.Range("E38").Value = "  -0.30%  "

# Row 39 - ARBITRUM (Volume only)
$ws.Range("E39").Value = "  +3.03%  "

# Row 40 - Pax Dollar (Volume only)
$ws.Range("E40").Value = "  -0.20%  "

# Row 41 - Frax Share
$ws.Range("D41").Value = "5.41"
$ws.Range("E41").Value = "  +2.45%  "

# Row 42 - MX Token (Volume only)
$ws.Range("E42").Value = "  +1.99%  "

# Row 43 - Trust Wallet Token
$ws.Range("D43").Value = "0.789"
$ws.Range("E43").Value = "  +0.41%  "

# Row 44 - Aave
$ws.Range("D44").Value = "63.83"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45 - Rocket Pool ETH
$ws.Range("D45").Value = "1.736.25"
$ws.Range("E45").Value = "  +0.86%  "

# Row 46 - WEMIX Token (Volume only)
$ws.Range("E46").Value = "  +7.35%  "

# Row 47 - Quant
$ws.Range("D47").Value = "90.17"
$ws.Range("E47").Value = "  +1.54%  "

# Row 48 - Render Token (Volume only)
$ws.Range("E48").Value = "  +0.38%  "

# Row 49 - Algorand (Volume only)
$ws.Range("E49").Value = "  +2.22%  "

# Row 50 - Cronos (Volume only)
$ws.Range("E50").Value = "  -0.38%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.51"
$ws.Range("E51").Value = "  +0.43%  "
